$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" '66.998.10'
Set-TextValue "E2" '  -0.56%  '
Set-TextValue "D3" '3.504.46'
Set-TextValue "E3" '  -2.31%  '
Set-TextValue "D4" '0.999'
Set-TextValue "E4" '  -0.09%  '
Set-TextValue "D5" '200.69'
Set-TextValue "E5" '  +4.01%  '
Set-TextValue "D6" '549.32'
Set-TextValue "E6" '  -4.97%  '
Set-TextValue "D7" '3.495.11'
Set-TextValue "E7" '  -2.42%  '
Set-TextValue "D8" '0.603'
Set-TextValue "E8" '  -2.37%  '
Set-TextValue "E9" '  +0.07%  '
Set-TextValue "D10" '0.655'
Set-TextValue "E10" '  -3.05%  '
Set-TextValue "D11" '60.70'
Set-TextValue "E11" '  +11.23%  '
Set-TextValue "E12" '  -4.54%  '
Set-TextValue "D13" '0.0000275'
Set-TextValue "E13" '  +0.32%  '
Set-TextValue "D14" '9.78'
Set-TextValue "E14" '  -1.51%  '
Set-TextValue "D15" '4.053.89'
Set-TextValue "E15" '  -2.45%  '
Set-TextValue "D16" '3.488.50'
Set-TextValue "E16" '  -3.09%  '
Set-TextValue "E17" '  -0.99%  '
Set-TextValue "D18" '18.38'
Set-TextValue "E18" '  -0.04%  '
Set-TextValue "D19" '66.620.48'
Set-TextValue "E19" '  -1.03%  '
Set-TextValue "D20" '11.80'
Set-TextValue "E20" '  -4.01%  '
Set-TextValue "E21" '  -3.78%  '
Set-TextValue "D22" '387.86'
Set-TextValue "E22" '  -3.51%  '
Set-TextValue "D23" '4.00'
Set-TextValue "E23" '  -4.87%  '
Set-TextValue "D24" '11.95'
Set-TextValue "E24" '  -11.25%  '
Set-TextValue "D25" '82.19'
Set-TextValue "E25" '  -3.99%  '
Set-TextValue "E26" '  +0.62%  '
Set-TextValue "D27" '2.80'
Set-TextValue "E27" '  -4.46%  '
Set-TextValue "D28" '11.89'
Set-TextValue "E28" '  -5.24%  '
Set-TextValue "E29" '  -2.50%  '
Set-TextValue "D30" '8.86'
Set-TextValue "E30" '  -3.06%  '
Set-TextValue "D31" '30.60'
Set-TextValue "E31" '  -2.21%  '
Set-TextValue "D32" '7.32'
Set-TextValue "E32" '  -9.32%  '
Set-TextValue "D33" '671.06'
Set-TextValue "E33" '  +0.31%  '
Set-TextValue "D34" '11.69'
Set-TextValue "E34" '  -4.21%  '
Set-TextValue "D35" '63.34'
Set-TextValue "E35" '  -1.00%  '
Set-TextValue "D36" '0.110'
Set-TextValue "E36" '  -5.49%  '
Set-TextValue "D37" '39.42'
Set-TextValue "E37" '  -7.24%  '
Set-TextValue "D38" '0.409'
Set-TextValue "E38" '  -3.71%  '
Set-TextValue "E39" '  -0.04%  '
Set-TextValue "D40" '3.06'
Set-TextValue "E40" '  -1.54%  '
Set-TextValue "D41" '3.111.76'
Set-TextValue "D42" '0.996'
Set-TextValue "E42" '  -0.20%  '
Set-TextValue "D43" '0.129'
Set-TextValue "E43" '  -3.62%  '
Set-TextValue "D44" '0.0₃0703'
Set-TextValue "E44" '  -10.37%  '
Set-TextValue "D45" '2.54'
Set-TextValue "E45" '  -13.27%  '
Set-TextValue "D46" '2.78'
Set-TextValue "E46" '  +15.14%  '
Set-TextValue "D47" '2.67'
Set-TextValue "E47" '  +6.37%  '
Set-TextValue "D48" '0.0397'
Set-TextValue "E48" '  -5.09%  '
Set-TextValue "E49" '  -3.60%  '
Set-TextValue "D50" '136.30'
Set-TextValue "E50" '  -4.27%  '
Set-TextValue "B51" 'ApeXProtocol'
Set-TextValue "C51" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue "D51" '2.93'
Set-TextValue "E51" '  -6.37%  '
